$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix marking scheme values (row 11) and totals (row 12) per corrected total marks
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -14
$ws.Range("E12").Value = "58 / 112"
